$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking price strings so they keep
# their original formatted text (e.g. "1.00") instead of becoming numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "66.813.79"
$ws.Range("E2").Value = "  +1.09%  "
$ws.Range("D3").Value = "3.097.72"
$ws.Range("E3").Value = "  +4.56%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "579.63"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").Value = "172.74"
$ws.Range("E6").Value = "  +3.71%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.091.90"
$ws.Range("E8").Value = "  +4.57%  "
$ws.Range("D9").Value = "0.523"
$ws.Range("E9").Value = "  +0.92%  "
$ws.Range("D10").Value = "6.43"
$ws.Range("E10").Value = "  -3.94%  "
$ws.Range("E11").Value = "  +2.21%  "
$ws.Range("E12").Value = "  +3.03%  "
$ws.Range("E13").Value = "  +1.09%  "
$ws.Range("D14").Value = "37.51"
$ws.Range("E14").Value = "  +5.26%  "
$ws.Range("E15").Value = "  -0.01%  "
$ws.Range("D16").Value = "3.609.63"
$ws.Range("E16").Value = "  +4.62%  "
$ws.Range("D17").Value = "66.789.08"
$ws.Range("E17").Value = "  +1.28%  "
$ws.Range("D18").Value = "7.20"
$ws.Range("E18").Value = "  +0.80%  "
$ws.Range("D19").Value = "3.097.44"
$ws.Range("E19").Value = "  +4.72%  "
$ws.Range("D20").Value = "16.27"
$ws.Range("E20").Value = "  +1.52%  "
$ws.Range("D21").Value = "480.68"
$ws.Range("E21").Value = "  +6.94%  "
$ws.Range("E22").Value = "  +2.12%  "
$ws.Range("E23").Value = "  +2.90%  "
$ws.Range("E24").Value = "  +2.00%  "
$ws.Range("D25").Value = "13.25"
$ws.Range("E25").Value = "  +6.79%  "
$ws.Range("E26").Value = "  +4.10%  "
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("D29").Value = "7.99"
$ws.Range("E29").Value = "  -3.48%  "
$ws.Range("E30").Value = "  -2.15%  "
$ws.Range("E31").Value = "  +2.80%  "
$ws.Range("D32").Value = "28.79"
$ws.Range("E32").Value = "  +4.78%  "
$ws.Range("E33").Value = "  -0.73%  "
$ws.Range("E34").Value = "  -2.49%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.29%  "
$ws.Range("E36").Value = "  +2.15%  "
$ws.Range("D37").Value = "0.990"
$ws.Range("E37").Value = "  +1.40%  "
$ws.Range("E38").Value = "  +0.68%  "
$ws.Range("E39").Value = "  +5.79%  "
$ws.Range("E40").Value = "  +3.80%  "
$ws.Range("D41").Value = "50.10"
$ws.Range("E41").Value = "  +1.79%  "
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("E43").Value = "  +1.75%  "
$ws.Range("E44").Value = "  -1.65%  "
$ws.Range("D45").Value = "2.841.87"
$ws.Range("E45").Value = "  +5.64%  "
$ws.Range("E46").Value = "  +1.98%  "
$ws.Range("D47").Value = "384.34"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("D48").Value = "135.53"
$ws.Range("E48").Value = "  +1.82%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").Value = "25.07"
$ws.Range("E50").Value = "  +3.94%  "
$ws.Range("E51").Value = "  +1.88%  "
